$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "made"
$ws.Range("B2").Value = "change"
$ws.Range("C2").Value = "to test"

$ws.Range("D6").Select()
